# Weekly refresh: a new week's Acelga price record is inserted at the top
# of this market's data block (row 489), pushing the existing rows
# (489:523) down by one row (to 490:524). The former last row (523)
# becomes row 524 with its data unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 489; rows 489:523 shift down to 490:524 and
# the sheet dimension grows to A1:R524 automatically.
$ws.Rows("489:489").Insert()

# Populate the newly inserted row 489 with this week's record.
$ws.Range("A489").Value = 3
$ws.Range("B489").Value = "Femacal de La Calera"
$ws.Range("C489").Value = "Coquimbo"
$ws.Range("D489").Value = 45021
$ws.Range("E489").Value = 5
$ws.Range("F489").Value = 100112009
$ws.Range("G489").Value = "Acelga"
$ws.Range("H489").Value = "Sin especificar"
$ws.Range("I489").Value = "Primera"
$ws.Range("J489").Value = 210
$ws.Range("K489").Value = 3500
$ws.Range("L489").Value = 3800
$ws.Range("M489").Value = 3657
$ws.Range("N489").Value = "`$/docena de atados (6 kilos)"
$ws.Range("O489").Value = "Provincia de Quillota"
$ws.Range("P489").Value = 610
$ws.Range("Q489").Value = 6
$ws.Range("R489").Value = "Hortaliza"
